# Add 2022-Q4 data to the workbook:
#  1. Insert a new row in the "总计" (total) summary sheet for 2022-Q4, pushing the
#     existing quarters down by one row.
#  2. Insert a brand-new worksheet named "2022-Q4" right before "2022-Q3" containing
#     the fund holding details for that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "总计" summary sheet (first sheet in the workbook)
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)

# Shift existing quarter rows (2..7) down by one row (8..3), working from the
# bottom up so we never overwrite data before it has been copied. Reading via
# .Formula (rather than .Value) is used since it reliably returns the actual
# cell content (text or number) in this runtime.
for ($r = 7; $r -ge 2; $r--) {
    $dst = $r + 1
    $totalSheet.Range("B$dst").Formula = $totalSheet.Range("B$r").Formula
    $totalSheet.Range("C$dst").Formula = $totalSheet.Range("C$r").Formula
    $totalSheet.Range("D$dst").Formula = $totalSheet.Range("D$r").Formula
}

# Give the newly-created row 8 the same look as the other index cells in column A.
$totalSheet.Range("A7").Copy()
$totalSheet.Range("A8").PasteSpecial(-4122)
$totalSheet.Range("A8").Value = 6

# Fill in the new 2022-Q4 figures at the top of the list (row 2).
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 0.53

# ---------------------------------------------------------------------------
# 2. Insert a new worksheet "2022-Q4" right before "2022-Q3"
# ---------------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q4Sheet = $wb.Worksheets.Add($q3Sheet)
$q4Sheet.Name = "2022-Q4"

# Copy header formatting (bold / centered / bordered) from the 2022-Q3 sheet.
$q3Sheet.Range("A1:H4").Copy()
$q4Sheet.Range("A1").PasteSpecial(-4122)

# Header row
$q4Sheet.Range("B1").Value = "基金代码"
$q4Sheet.Range("C1").Value = "基金名称"
$q4Sheet.Range("D1").Value = "基金规模"
$q4Sheet.Range("E1").Value = "股票总仓位"
$q4Sheet.Range("F1").Value = "仓位占比"
$q4Sheet.Range("G1").Value = "持有市值(亿元)"
$q4Sheet.Range("H1").Value = "仓位排名"

# Row 2
$q4Sheet.Range("A2").Value = 0
$q4Sheet.Range("B2").Value = "'002938"
$q4Sheet.Range("C2").Value = "中银证券健康产业灵活配置混合"
$q4Sheet.Range("D2").Value = "'5.33"
$q4Sheet.Range("E2").Value = "'92.95"
$q4Sheet.Range("F2").Value = "'4.76"
$q4Sheet.Range("G2").Value = "'0.2537"
$q4Sheet.Range("H2").Value = 2

# Row 3
$q4Sheet.Range("A3").Value = 1
$q4Sheet.Range("B3").Value = "'010054"
$q4Sheet.Range("C3").Value = "万家健康产业混合A"
$q4Sheet.Range("D3").Value = "'5.08"
$q4Sheet.Range("E3").Value = "'91.05"
$q4Sheet.Range("F3").Value = "'3.13"
$q4Sheet.Range("G3").Value = "'0.1590"
$q4Sheet.Range("H3").Value = 9

# Row 4
$q4Sheet.Range("A4").Value = 2
$q4Sheet.Range("B4").Value = "'010055"
$q4Sheet.Range("C4").Value = "万家健康产业混合C"
$q4Sheet.Range("D4").Value = "'3.85"
$q4Sheet.Range("E4").Value = "'91.05"
$q4Sheet.Range("F4").Value = "'3.13"
$q4Sheet.Range("G4").Value = "'0.1205"
$q4Sheet.Range("H4").Value = 9
